$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (Post 50)
$ws.Range("B60").Value = 50
$ws.Range("C60").Value = "Arithmetic Operations | Shell Scripting"
$ws.Range("D60").Value = 44174
$ws.Range("E60").Value = "https://programmingport.hashnode.dev/arithmetic-operations-or-shell-scripting"
$ws.Range("F60").Value = "https://dev.to/rahulmishra05/arithmetic-operations-shell-scripting-340j"

# Copy style formatting from the row above (row 59) so the new row matches
$ws.Range("B59:F59").Copy()
$ws.Range("B60:F60").PasteSpecial(-4122)

$ws.Range("B60").Value = 50
$ws.Range("C60").Value = "Arithmetic Operations | Shell Scripting"
$ws.Range("D60").Value = 44174
$ws.Range("E60").Value = "https://programmingport.hashnode.dev/arithmetic-operations-or-shell-scripting"
$ws.Range("F60").Value = "https://dev.to/rahulmishra05/arithmetic-operations-shell-scripting-340j"

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("B10:F60"))

# Update the sheet view (scroll position + selection) to match the commit
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("G61").Select()
